$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous template placeholder cells (A3:A5)
$ws.Range("A3:A5").Value = $null

# Body lines of the email get written out first, in reading order,
# to column A (skipping rows 3 and 5 as blank-line separators) and
# mirrored into column B.
$ws.Range("A2").Value = "This is a heading. "
$ws.Range("B2").Value = "This is a heading. "

$ws.Range("A4").Value = "More details. Test # 1. "
$ws.Range("B4").Value = "More details. Test # 1. "

$ws.Range("A6").Value = "Test Field #1: 204.33"
$ws.Range("B6").Value = "Test Field #1: 204.33"

$ws.Range("A7").Value = "Test Field #2: 201.23231"
$ws.Range("B7").Value = "Test Field #2: 201.23231"

$ws.Range("A8").Value = "Test Field #3: 701.9"
$ws.Range("B8").Value = "Test Field #3: 701.9"

# Finally, the message subject/id is written to the top row last.
$ws.Range("A1").Value = "#22222222222222"
$ws.Range("B1").Value = "#22222222222222"

$ws.Range("A1:F8").Select() | Out-Null
